$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("H68").Value = 39459
$ws.Range("J68").Value = 39459
$ws.Range("L68").Value = 39459
$ws.Range("N68").Value = -40957
$ws.Range("H71").Value = 39459
$ws.Range("J71").Value = 39459
$ws.Range("L71").Value = 118377
$ws.Range("N71").Value = -125865
$ws.Range("H132").Value = 4256982.5
$ws.Range("I132").Value = 4763482.5
$ws.Range("J132").Value = 2380
$ws.Range("K132").Value = 14290447.5
$ws.Range("L132").Value = 7140
$ws.Range("M132").Value = -14287917.5
$ws.Range("N132").Value = -12200
$ws.Range("H138").Value = 2238.1428
$ws.Range("I138").Value = 989.3774
$ws.Range("J138").Value = 4373.129
$ws.Range("K138").Value = 2968.1322
$ws.Range("L138").Value = 13119.387
$ws.Range("M138").Value = 2171.8678
$ws.Range("N138").Value = -23399.387

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4168.086
$ws.Range("I32").Value = 3331.1604
$ws.Range("J32").Value = 9817.333000000001
$ws.Range("K32").Value = 3331.1604
$ws.Range("L32").Value = 9817.333000000001
$ws.Range("M32").Value = -3044.1604
$ws.Range("N32").Value = -10391.333
$ws.Range("H74").Value = 657.63635
$ws.Range("I74").Value = 657.63635
$ws.Range("K74").Value = 657.63635
$ws.Range("M74").Value = 216.36365
$ws.Range("H77").Value = 657.63635
$ws.Range("I77").Value = 657.63635
$ws.Range("K77").Value = 3288.18175
$ws.Range("M77").Value = 1079.81825

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H99").Value = 4497.5
$ws.Range("I99").Value = 1995
$ws.Range("K99").Value = 1995
$ws.Range("M99").Value = -497
$ws.Range("H105").Value = 1983.5385
$ws.Range("I105").Value = 1606.1538
$ws.Range("J105").Value = 2360.923
$ws.Range("K105").Value = 1606.1538
$ws.Range("L105").Value = 2360.923
$ws.Range("M105").Value = 140.8462
$ws.Range("N105").Value = -5854.923

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 758.32355
$ws.Range("I5").Value = 488.1613
$ws.Range("J5").Value = 3550
$ws.Range("K5").Value = 1464.4839
$ws.Range("L5").Value = 10650
$ws.Range("M5").Value = -1352.4839
$ws.Range("N5").Value = -10874
$ws.Range("H20").Value = 2800
$ws.Range("J20").Value = 2750
$ws.Range("L20").Value = 8250
$ws.Range("N20").Value = -8704
$ws.Range("H40").Value = 116.666664
$ws.Range("I40").Value = 106.25
$ws.Range("J40").Value = 200
$ws.Range("K40").Value = 425
$ws.Range("L40").Value = 800
$ws.Range("M40").Value = -356
$ws.Range("N40").Value = -938
$ws.Range("H69").Value = 139612.5
$ws.Range("I69").Value = 1000
$ws.Range("J69").Value = 159414.28
$ws.Range("K69").Value = 3000
$ws.Range("L69").Value = 478242.84
$ws.Range("M69").Value = -2189
$ws.Range("N69").Value = -479864.84
$ws.Range("H72").Value = 139612.5
$ws.Range("I72").Value = 1000
$ws.Range("J72").Value = 159414.28
$ws.Range("K72").Value = 9000
$ws.Range("L72").Value = 1434728.52
$ws.Range("M72").Value = -4944
$ws.Range("N72").Value = -1442840.52
$ws.Range("H88").Value = 4465.077
$ws.Range("I88").Value = 2014
$ws.Range("J88").Value = 4669.3335
$ws.Range("K88").Value = 6042
$ws.Range("L88").Value = 14008.0005
$ws.Range("M88").Value = -5614
$ws.Range("N88").Value = -14864.0005
$ws.Range("H91").Value = 4465.077
$ws.Range("I91").Value = 2014
$ws.Range("J91").Value = 4669.3335
$ws.Range("K91").Value = 6042
$ws.Range("L91").Value = 14008.0005
$ws.Range("M91").Value = -4560
$ws.Range("N91").Value = -16972.0005
$ws.Range("H94").Value = 2950.9092
$ws.Range("J94").Value = 3307.5
$ws.Range("L94").Value = 9922.5
$ws.Range("N94").Value = -11274.5
$ws.Range("H104").Value = 2530
$ws.Range("I104").Value = 2200
$ws.Range("J104").Value = 2695
$ws.Range("K104").Value = 6600
$ws.Range("L104").Value = 8085
$ws.Range("M104").Value = -3979
$ws.Range("N104").Value = -13327
$ws.Range("H135").Value = 758.32355
$ws.Range("I135").Value = 488.1613
$ws.Range("J135").Value = 3550
$ws.Range("K135").Value = 4393.4517
$ws.Range("L135").Value = 31950
$ws.Range("M135").Value = -1858.4517
$ws.Range("N135").Value = -37020

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 4833.3335
$ws.Range("J58").Value = 4833.3335
$ws.Range("L58").Value = 4833.3335
$ws.Range("N58").Value = -5387.3335
$ws.Range("H62").Value = 10000
$ws.Range("J62").Value = 10000
$ws.Range("L62").Value = 10000
$ws.Range("N62").Value = -11372
$ws.Range("H65").Value = 10000
$ws.Range("J65").Value = 10000
$ws.Range("L65").Value = 30000
$ws.Range("N65").Value = -36864
$ws.Range("H70").Value = 4696.316
$ws.Range("I70").Value = 4884.2856
$ws.Range("J70").Value = 4170
$ws.Range("K70").Value = 4884.2856
$ws.Range("L70").Value = 4170
$ws.Range("M70").Value = -4614.2856
$ws.Range("N70").Value = -4710
$ws.Range("H73").Value = 4696.316
$ws.Range("I73").Value = 4884.2856
$ws.Range("J73").Value = 4170
$ws.Range("K73").Value = 4884.2856
$ws.Range("L73").Value = 4170
$ws.Range("M73").Value = -3948.2856
$ws.Range("N73").Value = -6042
$ws.Range("H122").Value = 3086.8
$ws.Range("I122").Value = 2380
$ws.Range("J122").Value = 3263.5
$ws.Range("K122").Value = 7140
$ws.Range("L122").Value = 9790.5
$ws.Range("M122").Value = -4690
$ws.Range("N122").Value = -14690.5

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2000
$ws.Range("J7").Value = 2000
$ws.Range("L7").Value = 2000
$ws.Range("N7").Value = -2224
$ws.Range("H126").Value = 2000
$ws.Range("J126").Value = 2000
$ws.Range("L126").Value = 6000
$ws.Range("N126").Value = -10940

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 70007
$ws.Range("J12").Value = 70007
$ws.Range("L12").Value = 70007
$ws.Range("N12").Value = -70291
$ws.Range("H81").Value = 1221.1111
$ws.Range("I81").Value = 998.3333
$ws.Range("J81").Value = 1666.6666
$ws.Range("K81").Value = 1996.6666
$ws.Range("L81").Value = 3333.3332
$ws.Range("M81").Value = -935.6666
$ws.Range("N81").Value = -5455.3332
$ws.Range("H84").Value = 1221.1111
$ws.Range("I84").Value = 998.3333
$ws.Range("J84").Value = 1666.6666
$ws.Range("K84").Value = 9983.333000000001
$ws.Range("L84").Value = 16666.666
$ws.Range("M84").Value = -4679.333000000001
$ws.Range("N84").Value = -27274.666
$ws.Range("H122").Value = 627007
$ws.Range("I122").Value = 715758
$ws.Range("K122").Value = 2147274
$ws.Range("M122").Value = -2144824
